# Add 2022-Q3 data:
#  - Insert a new worksheet "2022-Q3" right after the "总计" (totals) sheet,
#    carrying the per-fund holdings for the new quarter.
#  - Update the "总计" summary sheet with a new row for 2022-Q3 and shift the
#    existing quarters down by one row.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert the new "2022-Q3" sheet right after "总计"
# ------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)
$q3 = $wb.Worksheets.Add($null, $totalSheet)
$q3.Name = "2022-Q3"

# Header row
$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

$headerRange = $q3.Range("B1:H1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

# Data rows: A = running index (number), H = rank (number),
# everything else keeps its original text representation.
$q3Data = @(
    @("0", "630008", "华商策略精选混合",         "4.88", "77.93", "4.98", "0.2430", "2"),
    @("1", "011851", "天弘先进制造混合A",        "1.92", "91.41", "4.50", "0.0864", "7"),
    @("2", "011852", "天弘先进制造混合C",        "0.67", "91.41", "4.50", "0.0302", "7"),
    @("3", "001744", "诺安进取回报灵活配置混合", "0.23", "82.31", "3.82", "0.0088", "10")
)

$r = 2
foreach ($row in $q3Data) {
    $q3.Range("A$r").Value = [int]$row[0]
    $q3.Range("A$r").Font.Bold = $true
    $q3.Range("A$r").Borders.LineStyle = 1
    $q3.Range("A$r").HorizontalAlignment = -4108
    $q3.Range("A$r").VerticalAlignment = -4160

    $q3.Range("B$r").Value = "'" + $row[1]
    $q3.Range("C$r").Value = "'" + $row[2]
    $q3.Range("D$r").Value = "'" + $row[3]
    $q3.Range("E$r").Value = "'" + $row[4]
    $q3.Range("F$r").Value = "'" + $row[5]
    $q3.Range("G$r").Value = "'" + $row[6]
    $q3.Range("H$r").Value = [int]$row[7]
    $r++
}

# ------------------------------------------------------------------
# 2. Update the "总计" summary sheet
# ------------------------------------------------------------------
$totalData = @(
    @("2022-Q3", 4, 0.37),
    @("2022-Q2", 3, 0.42),
    @("2022-Q1", 4, 0.62),
    @("2021-Q4", 5, 0.91),
    @("2021-Q3", 4, 0.76),
    @("2021-Q2", 6, 0.87)
)

$r = 2
foreach ($row in $totalData) {
    $totalSheet.Range("A$r").Value = $r - 2
    $totalSheet.Range("B$r").Value = "'" + $row[0]
    $totalSheet.Range("C$r").Value = $row[1]
    $totalSheet.Range("D$r").Value = $row[2]
    $r++
}
